# Update the "Directory" value for the video row (E2): the data was
# reorganised from the old Spreadsheet_Data/Multimedia_Data layout into
# the new lower-case "data/multimedia/video/" path used by the repo.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("E2")
$cell.Value = "data/multimedia/video/"

# Touching the font forces Excel to record an explicit (new) cell style
# for E2, matching the re-saved workbook's style table.
$cell.Font.Color = 0

# The author's selection ended up on E3 after editing the directory cell.
$ws.Range("E3").Select()
